# Fix typo'd / concatenated header labels (missing spaces) across the
# various "Table_N" sheets, and fix the Table_6 (T-slot stud) table which
# had its data rows shifted by one (an extra header-only row at the top).

$wb = $excel.ActiveWorkbook

# --- Table_1 --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Table_1")
$ws.Range("D3").Value = "Mounting Hole Location"
$ws.Range("I3").Value = "Pkg. Qty."

# --- Table_2 --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Table_2")
$ws.Range("D3").Value = "Hole Location"
$ws.Range("E3").Value = "Fastener Thread Size"
$ws.Range("F3").Value = "Fastener Thread Lg."
$ws.Range("K3").Value = "Pkg. Qty."

# --- Table_3 --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Table_3")
$ws.Range("D3").Value = "Mounting Hole Location"
$ws.Range("I3").Value = "Pkg. Qty."

# --- Table_4 --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Table_4")
$ws.Range("B3").Value = "Double and Quad"
$ws.Range("D3").Value = "Hole Location"
$ws.Range("E3").Value = "For Screw Size"
$ws.Range("F3").Value = "Fastener Thread Size"
$ws.Range("G3").Value = "Fastener Thread Lg."
$ws.Range("L3").Value = "Pkg. Qty."

# --- Table_5 --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Table_5")
$ws.Range("C3").Value = "For Rail Ht., mm"
$ws.Range("D3").Value = "Thread Size"
$ws.Range("E3").Value = "Lg., mm"
$ws.Range("F3").Value = "Thick., mm"
$ws.Range("G3").Value = "Pkg. Qty."

# --- Table_6 ----------------------------------------------------------------
# Row 3 header label had a missing space.
# Row 4 was a stray header-only row (just "M8 x 1.25mm" thread size, with the
# rest of the row blank) that pushed the four real data rows (20mm/35mm/53mm/
# 71mm) down into rows 5-8, leaving row 8 an extra/duplicated trailing row.
# Deleting row 4 shifts the real data up into rows 4-7, and we restore the
# "M8 x 1.25mm" thread size into column K for each of those data rows.
$ws = $wb.Worksheets.Item("Table_6")
$ws.Range("B3").Value = "For T-Slot Wd."
$ws.Rows.Item(4).Delete()
$ws.Range("K4").Value = "M8 x 1.25mm"
$ws.Range("K5").Value = "M8 x 1.25mm"
$ws.Range("K6").Value = "M8 x 1.25mm"
$ws.Range("K7").Value = "M8 x 1.25mm"
